$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 40 with the daily auto-push data (2025/09/30, 火, 16, 12).
# Column A's value looks like a date ("2025/09/30"), so Excel would normally
# auto-convert it into a date serial number with a date number format.
# Force the cell to Text first so the value is kept as a literal string,
# then clear the leftover number-format override so no style index is
# left behind on the cell (matching the rest of the data rows, which are
# unstyled).
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = "2025/09/30"
$ws.Cells.Item(40, 1).ClearFormats()

$ws.Cells.Item(40, 2).Value = "火"
$ws.Cells.Item(40, 3).Value = 16
$ws.Cells.Item(40, 4).Value = 12
